$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Forecasted Consumption (MW)" values for A2:A97
$aValues = @(
    5550, 5500, 5460, 5420, 5380, 5360, 5350, 5340, 5330, 5320, 5320, 5320, 5320, 5320, 5330, 5350, 5380, 5410, 5440, 5480, 5520, 5560, 5610, 5660, 5720, 5780, 5820, 5880, 5940, 5990, 6020, 6050, 6060, 6060, 6050, 6050, 6030, 6010, 5990, 5970, 5950, 5950, 5950, 5950, 5950, 5950, 5950, 5960, 5980, 6000, 6000, 6000, 6000, 6000, 6000, 6020, 6080, 6150, 6210, 6280, 6340, 6440, 6550, 6680, 6840, 6950, 7030, 7090, 7100, 7100, 7100, 7100, 7060, 7040, 7020, 7000, 6950, 6900, 6830, 6760, 6660, 6550, 6450, 6330, 6190, 6080, 5930, 5820, 5720, 5630, 5540, 5440, 5390, 5350, 5280, 5220
)

# New "Timestamp" values for B2:B97 (Excel serial dates, +15 days vs previous run)
$bValues = @(
    45969, 45969.01041666666, 45969.02083333334, 45969.03125, 45969.04166666666, 45969.05208333334, 45969.0625, 45969.07291666666, 45969.08333333334, 45969.09375, 45969.10416666666, 45969.11458333334, 45969.125, 45969.13541666666, 45969.14583333334, 45969.15625, 45969.16666666666, 45969.17708333334, 45969.1875, 45969.19791666666, 45969.20833333334, 45969.21875, 45969.22916666666, 45969.23958333334, 45969.25, 45969.26041666666, 45969.27083333334, 45969.28125, 45969.29166666666, 45969.30208333334, 45969.3125, 45969.32291666666, 45969.33333333334, 45969.34375, 45969.35416666666, 45969.36458333334, 45969.375, 45969.38541666666, 45969.39583333334, 45969.40625, 45969.41666666666, 45969.42708333334, 45969.4375, 45969.44791666666, 45969.45833333334, 45969.46875, 45969.47916666666, 45969.48958333334, 45969.5, 45969.51041666666, 45969.52083333334, 45969.53125, 45969.54166666666, 45969.55208333334, 45969.5625, 45969.57291666666, 45969.58333333334, 45969.59375, 45969.60416666666, 45969.61458333334, 45969.625, 45969.63541666666, 45969.64583333334, 45969.65625, 45969.66666666666, 45969.67708333334, 45969.6875, 45969.69791666666, 45969.70833333334, 45969.71875, 45969.72916666666, 45969.73958333334, 45969.75, 45969.76041666666, 45969.77083333334, 45969.78125, 45969.79166666666, 45969.80208333334, 45969.8125, 45969.82291666666, 45969.83333333334, 45969.84375, 45969.85416666666, 45969.86458333334, 45969.875, 45969.88541666666, 45969.89583333334, 45969.90625, 45969.91666666666, 45969.92708333334, 45969.9375, 45969.94791666666, 45969.95833333334, 45969.96875, 45969.97916666666, 45969.98958333334
)

$rowCount = 96
$aArr = New-Object "object[,]" $rowCount,1
$bArr = New-Object "object[,]" $rowCount,1
for ($i = 0; $i -lt $rowCount; $i++) {
    $aArr[$i,0] = $aValues[$i]
    $bArr[$i,0] = $bValues[$i]
}

$ws.Range("A2:A97").Value = $aArr
$ws.Range("B2:B97").Value = $bArr

Write-Output "Updated $rowCount rows (A2:A97, B2:B97)"
